$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 3.45
$ws.Range("G2").Value = 4.2
$ws.Range("H2").Value = 2.06
$ws.Range("I2").Value = 2.32
$ws.Range("J2").Value = 3.4
$ws.Range("K2").Value = 3.95
$ws.Range("L2").Value = 1.39
$ws.Range("M2").Value = 1.06
$ws.Range("N2").Value = 4
$ws.Range("O2").Value = 1.28
$ws.Range("P2").Value = 1.98
$ws.Range("Q2").Value = 1.92
$ws.Range("R2").Value = 1.38
$ws.Range("S2").Value = 3.15
$ws.Range("T2").Value = 1.65
$ws.Range("U2").Value = 2.2
$ws.Range("V2").Value = 1.75
$ws.Range("W2").Value = 1.33
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 40
$ws.Range("AA2").Value = 30
$ws.Range("AB2").Value = 17.5
$ws.Range("AC2").Value = 8.6
$ws.Range("AD2").Value = 11.5
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 34
$ws.Range("AG2").Value = 16
$ws.Range("AH2").Value = 18.5
$ws.Range("AJ2").Value = 190
$ws.Range("AK2").Value = 44
$ws.Range("AL2").Value = 55
$ws.Range("AM2").Value = 85
$ws.Range("AN2").Value = 36
$ws.Range("AO2").Value = 16.5

# Row 3
$ws.Range("F3").Value = 6.4
$ws.Range("G3").Value = 8.199999999999999
$ws.Range("I3").Value = 1.49
$ws.Range("J3").Value = 5
$ws.Range("K3").Value = 6.4
$ws.Range("L3").Value = 1.25
$ws.Range("N3").Value = 6.2
$ws.Range("O3").Value = 1.16
$ws.Range("P3").Value = 2.84
$ws.Range("Q3").Value = 1.46
$ws.Range("R3").Value = 1.74
$ws.Range("S3").Value = 2.14
$ws.Range("T3").Value = 1.62
$ws.Range("U3").Value = 2.2
$ws.Range("W3").Value = 1.14
$ws.Range("Y3").Value = 1000
$ws.Range("Z3").Value = 1000
$ws.Range("AA3").Value = 900
$ws.Range("AD3").Value = 20
$ws.Range("AE3").Value = 1000
$ws.Range("AO3").Value = 5

# Row 4
$ws.Range("L4").Value = 1.39
$ws.Range("N4").Value = 3.4
$ws.Range("P4").Value = 1.86
$ws.Range("T4").Value = 1.87
$ws.Range("AG4").Value = 25

# Row 5
$ws.Range("F5").Value = 1.78
$ws.Range("G5").Value = 1.95
$ws.Range("H5").Value = 4.5
$ws.Range("I5").Value = 5.8
$ws.Range("J5").Value = 3.35
$ws.Range("K5").Value = 4.1
$ws.Range("L5").Value = 1.45
$ws.Range("N5").Value = 3.2
$ws.Range("Q5").Value = 2.08
$ws.Range("R5").Value = 1.27
$ws.Range("S5").Value = 3.75
$ws.Range("U5").Value = 1.85
$ws.Range("V5").Value = 1.21
$ws.Range("AB5").Value = 46
$ws.Range("AF5").Value = 55
$ws.Range("AG5").Value = 42

# Row 6
$ws.Range("H6").Value = 3.1
$ws.Range("J6").Value = 3.2
$ws.Range("L6").Value = 1.41
$ws.Range("N6").Value = 3.4
$ws.Range("P6").Value = 1.81
$ws.Range("Q6").Value = 1.93
$ws.Range("R6").Value = 1.31
$ws.Range("S6").Value = 3.35
$ws.Range("T6").Value = 1.72
$ws.Range("U6").Value = 2
$ws.Range("V6").Value = 1.35

# Row 7
$ws.Range("J7").Value = 3.75
$ws.Range("L7").Value = 1.37
$ws.Range("N7").Value = 3.45
$ws.Range("O7").Value = 1.29

# Row 8
$ws.Range("F8").Value = 2.06
$ws.Range("G8").Value = 2.12
$ws.Range("H8").Value = 4.2
$ws.Range("I8").Value = 5
$ws.Range("J8").Value = 3.2
$ws.Range("K8").Value = 3.65
$ws.Range("L8").Value = 1.51
$ws.Range("M8").Value = 1.1
$ws.Range("N8").Value = 2.98
$ws.Range("O8").Value = 1.43
$ws.Range("Q8").Value = 2.34
$ws.Range("R8").Value = 1.23
$ws.Range("S8").Value = 4.4
$ws.Range("T8").Value = 1.98
$ws.Range("U8").Value = 1.83
$ws.Range("V8").Value = 1.27
$ws.Range("W8").Value = 1.89
$ws.Range("X8").Value = 11.5
$ws.Range("Z8").Value = 36
$ws.Range("AA8").Value = 900
$ws.Range("AB8").Value = 7.6
$ws.Range("AC8").Value = 8.199999999999999
$ws.Range("AD8").Value = 23
$ws.Range("AF8").Value = 12
$ws.Range("AG8").Value = 11
$ws.Range("AH8").Value = 24
$ws.Range("AJ8").Value = 32
$ws.Range("AK8").Value = 27
$ws.Range("AL8").Value = 130
$ws.Range("AN8").Value = 23
$ws.Range("AO8").Value = 1000

# Row 9
$ws.Range("F9").Value = 2.42
$ws.Range("G9").Value = 2.6
$ws.Range("H9").Value = 3
$ws.Range("I9").Value = 3.3
$ws.Range("J9").Value = 3.35
$ws.Range("L9").Value = 1.41
$ws.Range("M9").Value = 1.07
$ws.Range("N9").Value = 3.65
$ws.Range("O9").Value = 1.32
$ws.Range("P9").Value = 1.9
$ws.Range("Q9").Value = 2.02
$ws.Range("S9").Value = 3.6
$ws.Range("T9").Value = 1.73
$ws.Range("U9").Value = 2.14
$ws.Range("V9").Value = 1.44
$ws.Range("W9").Value = 1.62
$ws.Range("X9").Value = 15.5
$ws.Range("Z9").Value = 44
$ws.Range("AA9").Value = 900
$ws.Range("AB9").Value = 11
$ws.Range("AC9").Value = 8.4
$ws.Range("AD9").Value = 24
$ws.Range("AE9").Value = 1000
$ws.Range("AG9").Value = 22
$ws.Range("AK9").Value = 1000
$ws.Range("AL9").Value = 1000
$ws.Range("AN9").Value = 980
$ws.Range("AO9").Value = 1000

# Row 10
$ws.Range("M10").Value = 1.02
$ws.Range("P10").Value = 1.51
$ws.Range("R10").Value = 1.51
$ws.Range("S10").Value = 1.05
$ws.Range("V10").Value = 1.01

# Row 11
$ws.Range("F11").Value = 1.59
$ws.Range("G11").Value = 1.61
$ws.Range("J11").Value = 4.9
$ws.Range("K11").Value = 5.1
$ws.Range("L11").Value = 1.27
$ws.Range("N11").Value = 6.8
$ws.Range("P11").Value = 2.96
$ws.Range("R11").Value = 1.79
$ws.Range("S11").Value = 2.22
$ws.Range("T11").Value = 1.6
$ws.Range("U11").Value = 2.62
$ws.Range("W11").Value = 2.62
$ws.Range("X11").Value = 32
$ws.Range("Y11").Value = 34
$ws.Range("Z11").Value = 130
$ws.Range("AF11").Value = 13
$ws.Range("AG11").Value = 11
$ws.Range("AH11").Value = 17
$ws.Range("AI11").Value = 55
$ws.Range("AK11").Value = 14
$ws.Range("AL11").Value = 23
$ws.Range("AM11").Value = 65
$ws.Range("AN11").Value = 5.7
$ws.Range("AO11").Value = 46

# Row 12
$ws.Range("F12").Value = 2.9
$ws.Range("G12").Value = 2.98
$ws.Range("H12").Value = 2.56
$ws.Range("I12").Value = 2.66
$ws.Range("J12").Value = 3.55
$ws.Range("K12").Value = 3.7
$ws.Range("L12").Value = 1.37
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 4.4
$ws.Range("O12").Value = 1.27
$ws.Range("P12").Value = 2.14
$ws.Range("Q12").Value = 1.85
$ws.Range("R12").Value = 1.45
$ws.Range("S12").Value = 3.05
$ws.Range("T12").Value = 1.67
$ws.Range("U12").Value = 2.34
$ws.Range("V12").Value = 1.6
$ws.Range("W12").Value = 1.5
$ws.Range("X12").Value = 17
$ws.Range("Y12").Value = 12.5
$ws.Range("Z12").Value = 18.5
$ws.Range("AB12").Value = 13.5
$ws.Range("AC12").Value = 8
$ws.Range("AE12").Value = 65
$ws.Range("AF12").Value = 21
$ws.Range("AG12").Value = 14
$ws.Range("AH12").Value = 15.5
$ws.Range("AI12").Value = 85
$ws.Range("AK12").Value = 38
$ws.Range("AL12").Value = 42
$ws.Range("AM12").Value = 330
$ws.Range("AN12").Value = 26
$ws.Range("AO12").Value = 44

# Row 13
$ws.Range("F13").Value = 3.4
$ws.Range("H13").Value = 2.08
$ws.Range("I13").Value = 2.12
$ws.Range("J13").Value = 4.1
$ws.Range("K13").Value = 4.3
$ws.Range("L13").Value = 1.3
$ws.Range("M13").Value = 1.03
$ws.Range("N13").Value = 5.8
$ws.Range("O13").Value = 1.2
$ws.Range("P13").Value = 2.62
$ws.Range("Q13").Value = 1.57
$ws.Range("R13").Value = 1.65
$ws.Range("S13").Value = 2.44
$ws.Range("T13").Value = 1.56
$ws.Range("U13").Value = 2.6
$ws.Range("V13").Value = 1.89
$ws.Range("X13").Value = 30
$ws.Range("Y13").Value = 14.5
$ws.Range("Z13").Value = 16
$ws.Range("AB13").Value = 38
$ws.Range("AC13").Value = 10.5
$ws.Range("AE13").Value = 19.5
$ws.Range("AF13").Value = 30
$ws.Range("AG13").Value = 15
$ws.Range("AI13").Value = 26
$ws.Range("AJ13").Value = 65
$ws.Range("AK13").Value = 34
$ws.Range("AL13").Value = 120
$ws.Range("AM13").Value = 60
$ws.Range("AN13").Value = 200
$ws.Range("AO13").Value = 9.800000000000001

# Row 14
$ws.Range("F14").Value = 6.6
$ws.Range("G14").Value = 7
$ws.Range("H14").Value = 1.58
$ws.Range("I14").Value = 1.59
$ws.Range("L14").Value = 1.35
$ws.Range("N14").Value = 4.7
$ws.Range("O14").Value = 1.25
$ws.Range("P14").Value = 2.24
$ws.Range("Q14").Value = 1.76
$ws.Range("R14").Value = 1.48
$ws.Range("S14").Value = 2.9
$ws.Range("T14").Value = 1.86
$ws.Range("V14").Value = 2.68
$ws.Range("W14").Value = 1.16
$ws.Range("X14").Value = 18.5
$ws.Range("Z14").Value = 9.6
$ws.Range("AA14").Value = 14.5
$ws.Range("AC14").Value = 9.800000000000001
$ws.Range("AG14").Value = 24
$ws.Range("AI14").Value = 32
$ws.Range("AJ14").Value = 190
$ws.Range("AL14").Value = 95
$ws.Range("AM14").Value = 110
$ws.Range("AN14").Value = 260
$ws.Range("AO14").Value = 7.6
